$d = $word.ActiveDocument

$d.Content.Find.Execute("2000", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2005", 2)

$d.Content.Find.Execute("Mengenorientierte Auswertung von Anfragen in der Logikprogrammiersprache PROLOG", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MetaObject Protocol Concepts for a RISC Object Model.", 2)

$d.Content.Find.Execute("Erich Gehlen, Burkhard Kehrbusch", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Frank Manola", 2)

$d.Content.Find.Execute("BSP Business School Berlin", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fachhochschule Polizei Sachsen-Anhalt (Aschersleben)", 2)

$d.Content.Find.Execute("Schmarjestrasse 32", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Prager Str 23", 2)
